$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before Q (shifts old Q:U -> R:V), carrying
#     values/styles/column-widths along the way, matching how the
#     "# Fish Released" field was inserted into the collections sheet.
$pw = $ws.Columns("P:P").ColumnWidth
$ws.Columns("Q:Q").Insert()
$ws.Columns("Q:Q").ColumnWidth = $pw

# New header cell for the inserted column.
$ws.Range("Q3").Value = "# Fish Released"

# --- Shift the cell-comments that lived on the old Q3:U3 header cells one
#     column to the right so they stay attached to the same logical field
#     (Temperature, Fishing Settings, Fishing Seconds, Voltage, Comments).
$txtQ = $ws.Range("Q3").Comment.Text()
$txtR = $ws.Range("R3").Comment.Text()
$txtS = $ws.Range("S3").Comment.Text()
$txtT = $ws.Range("T3").Comment.Text()
$txtU = $ws.Range("U3").Comment.Text()

$ws.Range("R3").Comment.Text($txtQ)
$ws.Range("S3").Comment.Text($txtR)
$ws.Range("T3").Comment.Text($txtS)
$ws.Range("U3").Comment.Text($txtT)

$ws.Range("Q3").Comment.Delete()
$ws.Range("V3").AddComment($txtU)

# --- Match the saved selection.
$ws.Range("Q4").Select()
